$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 341 (shifts existing rows 341-364 down to 342-365),
# matching the weekly data refresh described in the commit message.
$ws.Rows.Item(341).Insert()

$ws.Cells.Item(341, 1).Value = 3
$ws.Cells.Item(341, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(341, 3).Value = "Coquimbo"
$ws.Cells.Item(341, 4).Value = 44746
$ws.Cells.Item(341, 5).Value = 5
$ws.Cells.Item(341, 6).Value = 100112009
$ws.Cells.Item(341, 7).Value = "Acelga"
$ws.Cells.Item(341, 8).Value = "Sin especificar"
$ws.Cells.Item(341, 9).Value = "Primera"
$ws.Cells.Item(341, 10).Value = 215
$ws.Cells.Item(341, 11).Value = 3000
$ws.Cells.Item(341, 12).Value = 3500
$ws.Cells.Item(341, 13).Value = 3267
$ws.Cells.Item(341, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(341, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(341, 16).Value = 544
$ws.Cells.Item(341, 17).Value = 6
$ws.Cells.Item(341, 18).Value = "Hortaliza"
